# Apply the cryptocurrency price/volume refresh described in the commit.
# Values that look like plain numbers (e.g. "217.40") are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inline-string cell contents) instead of silently coercing them
# to numeric values and dropping formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.052.78"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.651.59"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'217.40"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'0.5259"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.2597"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "'0.06322"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'20.35"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").Value = "'0.07794"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "'4.503"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "1.643.56"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'0.5485"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "0.0₅8189"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "'65.51"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "26.069.15"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "'191.02"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "'10.09"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'6.029"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'143.42"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "'0.1236"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'7.216"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'0.05804"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "'1.273"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'3.554"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'3.268"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'1.581"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'2.780"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'0.9453"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.409"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "'0.5734"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "'0.01609"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.8425"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.741"
$ws.Range("E40").Value = "  -5.08%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'103.62"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").Value = "1.028.82"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "1.795.42"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'56.81"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'0.4323"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").Value = "'7.846"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "'0.05143"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "'1.465"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.09623"
$ws.Range("E51").Value = "  -0.73%  "
